$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
# Force text format on the Price cells before writing so Excel COM
# does not auto-coerce numeric-looking strings (e.g. "291.66") into
# floating point numbers; the source data stores these as literal text.
$dPriceCells = @("D2","D3","D6","D7","D8","D9","D13","D15","D16","D17","D18","D19","D21","D22","D23","D24","D25","D26","D27","D28","D30","D31","D32","D33","D34","D35","D36","D37","D39","D40","D41","D42","D43","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '22.438.02'
$ws.Range("D3").Value = '1.572.90'
$ws.Range("D6").Value = '291.66'
$ws.Range("D7").Value = '0.3734'
$ws.Range("D8").Value = '49.81'
$ws.Range("D9").Value = '0.3392'
$ws.Range("D13").Value = '21.38'
$ws.Range("D15").Value = '6.929'
$ws.Range("D16").Value = '1.570.50'
$ws.Range("D17").Value = '0.00001121'
$ws.Range("D18").Value = '91.00'
$ws.Range("D19").Value = '0.06735'
$ws.Range("D21").Value = '6.284'
$ws.Range("D22").Value = '16.36'
$ws.Range("D23").Value = '12.15'
$ws.Range("D24").Value = '22.433.84'
$ws.Range("D25").Value = '2.334'
$ws.Range("D26").Value = '2.625'
$ws.Range("D27").Value = '20.09'
$ws.Range("D28").Value = '148.38'
$ws.Range("D30").Value = '125.63'
$ws.Range("D31").Value = '1.747.53'
$ws.Range("D32").Value = '1.049'
$ws.Range("D33").Value = '6.134'
$ws.Range("D34").Value = '1.980'
$ws.Range("D35").Value = '9.778'
$ws.Range("D36").Value = '0.08358'
$ws.Range("D37").Value = '1.385'
$ws.Range("D39").Value = '0.2286'
$ws.Range("D40").Value = '0.06514'
$ws.Range("D41").Value = '5.459'
$ws.Range("D42").Value = '11.29'
$ws.Range("D43").Value = '0.6217'
$ws.Range("D45").Value = '13.86'
$ws.Range("D46").Value = '3.812'
$ws.Range("D47").Value = '0.5804'
$ws.Range("D48").Value = '129.62'
$ws.Range("D49").Value = '2.076'
$ws.Range("D50").Value = '1.216'
$ws.Range("D51").Value = '0.07322'

foreach ($addr in $dPriceCells) {
    $ws.Range($addr).Style = "Normal"
}

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("E6").Value = '  +0.30%  '
$ws.Range("E7").Value = '  -0.91%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  -0.84%  '
$ws.Range("E10").Value = '  -1.27%  '
$ws.Range("E11").Value = '  -2.13%  '
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("E13").Value = '  +0.62%  '
$ws.Range("E14").Value = '  -0.45%  '
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("E17").Value = '  -1.23%  '
$ws.Range("E18").Value = '  +1.13%  '
$ws.Range("E19").Value = '  -0.46%  '
$ws.Range("E21").Value = '  +1.14%  '
$ws.Range("E22").Value = '  -2.56%  '
$ws.Range("E23").Value = '  +0.97%  '
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("E25").Value = '  -3.58%  '
$ws.Range("E26").Value = '  -3.75%  '
$ws.Range("E27").Value = '  -1.01%  '
$ws.Range("E28").Value = '  +1.21%  '
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("E30").Value = '  -0.47%  '
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("E32").Value = '  +5.74%  '
$ws.Range("E33").Value = '  -0.83%  '
$ws.Range("E34").Value = '  -1.50%  '
$ws.Range("E35").Value = '  -2.46%  '
$ws.Range("E36").Value = '  -2.79%  '
$ws.Range("E37").Value = '  +4.19%  '
$ws.Range("E38").Value = '  -3.29%  '
$ws.Range("E39").Value = '  -1.16%  '
$ws.Range("E40").Value = '  -0.84%  '
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("E42").Value = '  -2.00%  '
$ws.Range("E43").Value = '  -3.72%  '
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("E45").Value = '  -1.96%  '
$ws.Range("E46").Value = '  +0.37%  '
$ws.Range("E47").Value = '  -3.44%  '
$ws.Range("E48").Value = '  +3.31%  '
$ws.Range("E49").Value = '  -0.39%  '
$ws.Range("E50").Value = '  -6.86%  '
$ws.Range("E51").Value = '  -0.14%  '

Write-Output "Applied cryptos price/volume update"
